$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '41.802.34'
$ws.Cells.Item(2, 5).Value = '  +4.12%  '

$ws.Cells.Item(3, 4).Value = '2.264.45'
$ws.Cells.Item(3, 5).Value = '  +1.89%  '

$ws.Cells.Item(4, 5).Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '305.30'
$ws.Cells.Item(5, 5).Value = '  +3.83%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '91.87'
$ws.Cells.Item(6, 5).Value = '  +4.55%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.532'
$ws.Cells.Item(7, 5).Value = '  +3.56%  '

$ws.Cells.Item(8, 5).Value = '  -0.02%  '

$ws.Cells.Item(9, 5).Value = '  +2.54%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '32.52'
$ws.Cells.Item(10, 5).Value = '  +6.02%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '53.92'
$ws.Cells.Item(11, 5).Value = '  +5.97%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.0796'
$ws.Cells.Item(12, 5).Value = '  +1.75%  '

$ws.Cells.Item(13, 5).Value = '  +0.89%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '6.60'
$ws.Cells.Item(14, 5).Value = '  +2.97%  '

$ws.Cells.Item(15, 4).Value = '2.615.25'
$ws.Cells.Item(15, 5).Value = '  +1.92%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '14.20'
$ws.Cells.Item(16, 5).Value = '  +2.64%  '

$ws.Cells.Item(17, 4).Value = '2.264.77'
$ws.Cells.Item(17, 5).Value = '  +1.85%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '0.762'
$ws.Cells.Item(18, 5).Value = '  +3.51%  '

$ws.Cells.Item(19, 4).Value = '41.704.89'
$ws.Cells.Item(19, 5).Value = '  +4.08%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '12.23'
$ws.Cells.Item(20, 5).Value = '  +8.50%  '

$ws.Cells.Item(21, 4).Value = '0.0₃0905'
$ws.Cells.Item(21, 5).Value = '  +1.72%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '5.91'
$ws.Cells.Item(22, 5).Value = '  +2.13%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '66.81'
$ws.Cells.Item(23, 5).Value = '  +1.75%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '240.90'
$ws.Cells.Item(24, 5).Value = '  +2.13%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '2.59'
$ws.Cells.Item(25, 5).Value = '  +4.75%  '

$ws.Cells.Item(26, 5).Value = '  +0.11%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '1.91'
$ws.Cells.Item(27, 5).Value = '  +5.15%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '24.14'
$ws.Cells.Item(28, 5).Value = '  +3.94%  '

$ws.Cells.Item(29, 5).Value = '  +11.29%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '9.58'
$ws.Cells.Item(30, 5).Value = '  +2.59%  '

$ws.Cells.Item(31, 2).Value = 'Monero'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D31").NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '158.17'
$ws.Cells.Item(31, 5).Value = '  -0.91%  '

$ws.Cells.Item(32, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D32").NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '33.83'
$ws.Cells.Item(32, 5).Value = '  +6.30%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '1.00'

$ws.Range("D34").NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '5.15'
$ws.Cells.Item(34, 5).Value = '  +3.69%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.0746'
$ws.Cells.Item(35, 5).Value = '  +4.25%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '3.00'
$ws.Cells.Item(36, 5).Value = '  -1.51%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '16.93'
$ws.Cells.Item(37, 5).Value = '  +8.19%  '

$ws.Cells.Item(38, 5).Value = '  +1.57%  '

$ws.Cells.Item(39, 5).Value = '  +2.47%  '

$ws.Cells.Item(40, 5).Value = '  +3.78%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '1.80'
$ws.Cells.Item(41, 5).Value = '  +1.76%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '3.91'
$ws.Cells.Item(42, 5).Value = '  +4.13%  '

$ws.Cells.Item(43, 4).Value = '2.062.05'
$ws.Cells.Item(43, 5).Value = '  -0.61%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '19.32'
$ws.Cells.Item(44, 5).Value = '  -1.35%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.0278'
$ws.Cells.Item(45, 5).Value = '  +2.68%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '10.25'
$ws.Cells.Item(46, 5).Value = '  +2.68%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '2.89'
$ws.Cells.Item(47, 5).Value = '  +4.31%  '

$ws.Cells.Item(48, 5).Value = '  +6.89%  '

$ws.Cells.Item(49, 5).Value = '  +4.03%  '

$ws.Cells.Item(50, 5).Value = '  +2.65%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '72.65'
$ws.Cells.Item(51, 5).Value = '  +7.01%  '
